# "final files from john"
# Add Donald J. Trump as the 45th President to the US Presidents table,
# close out Barack Obama's term (Date left Office), bump George W. Bush's
# "Date took office" cell to the new MM/DD/YYYY number format, and update
# the workbook/sheet view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Obama (row 45): he left office 2017-01-20 -------------------------
$ws.Range("I45").NumberFormat = "@"
$ws.Range("I45").Value = "2017-01-20"
$ws.Range("I45").NumberFormat = $ws.Range("H45").NumberFormat

# --- George W. Bush (row 44): reformat "Date took office" --------------
# Same text value (2001-01-20), just a new number format (MM/DD/YYYY)
$ws.Range("H44").NumberFormat = "MM/DD/YYYY"

# --- New row 46: Donald J. Trump, 45th president ------------------------
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "Trump"
$ws.Range("C46").Value = "Donald J"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1946-06-14"
$ws.Range("D46").NumberFormat = $ws.Range("D45").NumberFormat

$ws.Range("E46").Value = "NONE"

$ws.Range("F46").Value = "Queens, NYC"
$ws.Range("G46").Value = "New York"

$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "2017-01-20"
$ws.Range("H46").NumberFormat = $ws.Range("H45").NumberFormat

$ws.Range("I46").Value = "NONE"
$ws.Range("J46").Value = "Republican"

# --- Workbook / sheet view bookkeeping ----------------------------------
try { $wb.TabRatio = 988 } catch {}
try { $excel.ActiveWindow.TabRatio = 988 } catch {}

$ws.Range("D47").Select() | Out-Null
